# Applies the "Remove Fire/Lightning/Energy damage type multipliers" edit to
# the "Menu Options" sheet:
#   - Row 9  (OptionSlashMultiplier): default value changes 1.0f -> 0.8f
#   - Row 10 (OptionBluntMultiplier): unchanged
#   - Row 11 (OptionFireMultiplier) is repurposed into OptionBurningMultiplier
#   - Row 12 (OptionLightningMultiplier) is repurposed into OptionElectrocuteMultiplier
#   - Row 13 (OptionEnergyMultiplier) is removed entirely, shifting every
#     subsequent row up by one (rows 14-57 become rows 13-56)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Slash multiplier default (row 9, column E).
$ws.Cells.Item(9, 5).Value = "0.8f"

# 2. Repurpose row 11 (Fire -> Burning).
$ws.Cells.Item(11, 2).Value = 70
$ws.Cells.Item(11, 3).Value = "OptionBurningMultiplier"
$ws.Cells.Item(11, 5).Value = "0.5f"
$ws.Cells.Item(11, 6).Value = "Additional DOT damage multiplier when creature has Burning status. Stacks with Fire multiplier. 0.0x = no bonus damage from burning status."

# 3. Repurpose row 12 (Lightning -> Electrocute).
$ws.Cells.Item(12, 2).Value = 80
$ws.Cells.Item(12, 3).Value = "OptionElectrocuteMultiplier"
$ws.Cells.Item(12, 5).Value = "0.8f"
$ws.Cells.Item(12, 6).Value = "DOT damage multiplier when creature has Electrocute status. Electrocute normally does no damage. 0.0x = no damage from electrocute."

# 4. Remove row 13 (OptionEnergyMultiplier) entirely. This shifts every
#    following row up by one, which naturally re-numbers all the Zone /
#    Advanced / Statistics rows to match the target layout (and shrinks the
#    used range from G57 to G56) without having to touch them individually.
$ws.Rows.Item(13).Delete()

$wb.Save()
